$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: student record (file no, name, subjects, totals) ---
$ws.Range("A6").Value = 15611
$ws.Range("B6").Value = "dan"
$ws.Range("C6").ClearContents()
$ws.Range("D6").Value = 87
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = 74
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("J6").Value = 76
$ws.Range("K6").Value = 237
$ws.Range("L6").Value = 29.625
$ws.Range("M6").Value = "E"

# --- Row 7: student record ---
$ws.Range("A7").Value = 15612
$ws.Range("B7").Value = "tes"
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 89
$ws.Range("E7").ClearContents()
$ws.Range("F7").Value = 48
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("J7").Value = 78
$ws.Range("K7").Value = 215
$ws.Range("L7").Value = 26.875
$ws.Range("M7").Value = "E"

# --- Row 8: SubjectTotal ---
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 176
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 122
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 154
$ws.Range("K8").Value = 452
$ws.Range("L8").Value = 56.5

# --- Row 9: SubjectAverage ---
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 88
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 61
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 77
$ws.Range("K9").Value = 226
$ws.Range("L9").Value = 28.25

# --- Row 10: SubjectGrades ---
$ws.Range("C10").Value = "E"
$ws.Range("D10").Value = "A-"
$ws.Range("E10").Value = "E"
$ws.Range("F10").Value = "C+"
$ws.Range("G10").Value = "E"
$ws.Range("H10").Value = "E"
$ws.Range("I10").Value = "E"
$ws.Range("J10").Value = "B+"
$ws.Range("L10").Value = "E"
